$wb = $excel.ActiveWorkbook

# Target "ideal" (autofit) column width from the authoring Excel session.
# The headless engine quantizes ColumnWidth to whole pixels (grid of 1/6),
# so we feed it the input that lands on the pixel bucket nearest the
# canonical target width (17.2159881591797 ~= 17.1667 after quantization).
$targetColWidth = 16.3333333333333

# --- Sheet "Overview" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-08-22 22:57:34"
$ws1.Columns.Item(5).ColumnWidth = $targetColWidth
$ws1.Columns.Item(6).ColumnWidth = $targetColWidth

# --- Sheet "zh-cn" ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-08-22 22:57:29"
$ws2.Columns.Item(3).ColumnWidth = $targetColWidth

# --- Sheet "de-de" ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-08-22 22:57:34"
$ws3.Columns.Item(3).ColumnWidth = $targetColWidth
